# feat: update task buff tower
# 更新部分塔 更新buff 更新任务奖励 更新世界3波次奖励
#
# Updates the "rewards" (column G) and related numeric (column I) values
# for several rows in the Task (任务) sheet, then leaves the selection
# on H16 as in the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23/24 reward strings changed.
$ws.Range("G23").Value = "1|400||2|200||3|10"
$ws.Range("G24").Value = "1|500||2|300||3|10"
$ws.Range("I24").Value = 180

# World-3 wave reward tuning (rows 46-48).
$ws.Range("I46").Value = 15
$ws.Range("I47").Value = 25
$ws.Range("I48").Value = 35

# Tower buff rewards (rows 50-52).
$ws.Range("G50").Value = "1|600||2|300"
$ws.Range("G51").Value = "1|800||2|400"
$ws.Range("G52").Value = "1|1000||2|500"
$ws.Range("I52").Value = 130

# Row 54 reward + amount update.
$ws.Range("G54").Value = "1|1200||2|600||3|60"
$ws.Range("I54").Value = 2000

# Rows 55/56 amount reduced.
$ws.Range("I55").Value = 200
$ws.Range("I56").Value = 200

# Rows 71/72/74 amount updates.
$ws.Range("I71").Value = 10
$ws.Range("I72").Value = 10
$ws.Range("I74").Value = 100

# Restore the author's on-save selection / scroll position.
$ws.Range("H16").Select()
